# Auto-generated edit script applying the diff to Aegis_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1027.0555
$ws.Range("J112").Value = 1036.6875
$ws.Range("L112").Value = 3110.0625
$ws.Range("N112").Value = -5326.0625
$ws.Range("H121").Value = 619.6129
$ws.Range("J121").Value = 615.93335
$ws.Range("L121").Value = 1847.80005
$ws.Range("N121").Value = -5341.80005
$ws.Range("H127").Value = 21278572
$ws.Range("I127").Value = 685.875
$ws.Range("J127").Value = 25643266
$ws.Range("K127").Value = 2057.625
$ws.Range("L127").Value = 76929798
$ws.Range("M127").Value = 2902.375
$ws.Range("N127").Value = -76939718
$ws.Range("H132").Value = 9623610
$ws.Range("I132").Value = 10008450
$ws.Range("K132").Value = 30025350
$ws.Range("M132").Value = -30022820
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H137").Value = 1350.1111
$ws.Range("I137").Value = 1074.84
$ws.Range("K137").Value = 3224.52
$ws.Range("M137").Value = -674.5199999999995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22930.52
$ws.Range("I32").Value = 3786.7231
$ws.Range("J32").Value = 126626.086
$ws.Range("K32").Value = 3786.7231
$ws.Range("L32").Value = 126626.086
$ws.Range("M32").Value = -3499.7231
$ws.Range("N32").Value = -127200.086
$ws.Range("H102").Value = 102586.9
$ws.Range("I102").Value = 252244.75
$ws.Range("J102").Value = 2815
$ws.Range("K102").Value = 252244.75
$ws.Range("L102").Value = 2815
$ws.Range("M102").Value = -250622.75
$ws.Range("N102").Value = -6059
$ws.Range("H132").Value = 4576.125
$ws.Range("I132").Value = 5457.44
$ws.Range("J132").Value = 1428.5714
$ws.Range("K132").Value = 16372.32
$ws.Range("L132").Value = 4285.7142
$ws.Range("M132").Value = -13842.32
$ws.Range("N132").Value = -9345.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 11533.333
$ws.Range("J61").Value = 11533.333
$ws.Range("L61").Value = 11533.333
$ws.Range("N61").Value = -12159.333
$ws.Range("H105").Value = 68455.7
$ws.Range("I105").Value = 43290.832
$ws.Range("J105").Value = 169115.17
$ws.Range("K105").Value = 43290.832
$ws.Range("L105").Value = 169115.17
$ws.Range("M105").Value = -41543.832
$ws.Range("N105").Value = -172609.17
$ws.Range("H107").Value = 14512887
$ws.Range("I107").Value = 17567172
$ws.Range("K107").Value = 17567172
$ws.Range("M107").Value = -17565252

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23218.086
$ws.Range("I31").Value = 53837.95
$ws.Range("J31").Value = 2440.3215
$ws.Range("K31").Value = 53837.95
$ws.Range("L31").Value = 2440.3215
$ws.Range("M31").Value = -53542.95
$ws.Range("N31").Value = -3030.3215
$ws.Range("H34").Value = 23218.086
$ws.Range("I34").Value = 53837.95
$ws.Range("J34").Value = 2440.3215
$ws.Range("K34").Value = 53837.95
$ws.Range("L34").Value = 2440.3215
$ws.Range("M34").Value = -53635.95
$ws.Range("N34").Value = -2844.3215
$ws.Range("H70").Value = 16500
$ws.Range("J70").Value = 16500
$ws.Range("L70").Value = 16500
$ws.Range("N70").Value = -17130
$ws.Range("H73").Value = 16500
$ws.Range("J73").Value = 16500
$ws.Range("L73").Value = 16500
$ws.Range("N73").Value = -18684
$ws.Range("H99").Value = 7492.6816
$ws.Range("I99").Value = 2794.8462
$ws.Range("J99").Value = 14278.444
$ws.Range("K99").Value = 2794.8462
$ws.Range("L99").Value = 14278.444
$ws.Range("M99").Value = -1296.8462
$ws.Range("N99").Value = -17274.444
$ws.Range("H126").Value = 7492.6816
$ws.Range("I126").Value = 2794.8462
$ws.Range("J126").Value = 14278.444
$ws.Range("K126").Value = 8384.5386
$ws.Range("L126").Value = 42835.33199999999
$ws.Range("M126").Value = -5914.5386
$ws.Range("N126").Value = -47775.33199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 500
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 1500
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -3996
$ws.Range("H94").Value = 3741.6667
$ws.Range("I94").Value = 2983.3333
$ws.Range("J94").Value = 4500
$ws.Range("K94").Value = 8949.999899999999
$ws.Range("L94").Value = 13500
$ws.Range("M94").Value = -8273.999899999999
$ws.Range("N94").Value = -14852
$ws.Range("H113").Value = 1196.9333
$ws.Range("I113").Value = 503.33334
$ws.Range("J113").Value = 1890.5333
$ws.Range("K113").Value = 1510.00002
$ws.Range("L113").Value = 5671.5999
$ws.Range("M113").Value = 659.9999800000001
$ws.Range("N113").Value = -10011.5999
$ws.Range("H131").Value = 921.88
$ws.Range("I131").Value = 690
$ws.Range("J131").Value = 934.0842
$ws.Range("K131").Value = 2070
$ws.Range("L131").Value = 2802.2526
$ws.Range("M131").Value = 2970
$ws.Range("N131").Value = -12882.2526

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H119").Value = 44370.5
$ws.Range("J119").Value = 44370.5
$ws.Range("L119").Value = 44370.5
$ws.Range("N119").Value = -54046.5
$ws.Range("H122").Value = 2900.45
$ws.Range("I122").Value = 2706.5
$ws.Range("K122").Value = 8119.5
$ws.Range("M122").Value = -5669.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1463.9678
$ws.Range("J7").Value = 1879.1538
$ws.Range("L7").Value = 1879.1538
$ws.Range("N7").Value = -2103.1538
$ws.Range("H46").Value = 461086.3
$ws.Range("I46").Value = 576.9231
$ws.Range("J46").Value = 1126266.5
$ws.Range("K46").Value = 576.9231
$ws.Range("L46").Value = 1126266.5
$ws.Range("M46").Value = -388.9231
$ws.Range("N46").Value = -1126642.5
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 5000
$ws.Range("K62").Value = 5000
$ws.Range("M62").Value = -4376
$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 5000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880
$ws.Range("H82").Value = 1514
$ws.Range("I82").Value = 1462.2632
$ws.Range("J82").Value = 1759.75
$ws.Range("K82").Value = 1462.2632
$ws.Range("L82").Value = 1759.75
$ws.Range("M82").Value = -1101.2632
$ws.Range("N82").Value = -2481.75
$ws.Range("H85").Value = 1514
$ws.Range("I85").Value = 1462.2632
$ws.Range("J85").Value = 1759.75
$ws.Range("K85").Value = 1462.2632
$ws.Range("L85").Value = 1759.75
$ws.Range("M85").Value = -214.2632000000001
$ws.Range("N85").Value = -4255.75
$ws.Range("H122").Value = 4001.6667
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 4002.5
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 12007.5
$ws.Range("M122").Value = -9550
$ws.Range("N122").Value = -16907.5
$ws.Range("H126").Value = 1463.9678
$ws.Range("J126").Value = 1879.1538
$ws.Range("L126").Value = 5637.4614
$ws.Range("N126").Value = -10577.4614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 286614
$ws.Range("I81").Value = 167716.5
$ws.Range("K81").Value = 335433
$ws.Range("M81").Value = -334372
$ws.Range("H84").Value = 286614
$ws.Range("I84").Value = 167716.5
$ws.Range("K84").Value = 1677165
$ws.Range("M84").Value = -1671861
$ws.Range("H122").Value = 1026.8572
$ws.Range("I122").Value = 1039.6
$ws.Range("K122").Value = 3118.8
$ws.Range("M122").Value = -668.7999999999997
$ws.Range("H136").Value = 842.675
$ws.Range("I136").Value = 594.2963
$ws.Range("J136").Value = 1358.5385
$ws.Range("K136").Value = 1782.8889
$ws.Range("L136").Value = 4075.6155
$ws.Range("M136").Value = 767.1111000000001
$ws.Range("N136").Value = -9175.6155
